$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.224.95'
$ws.Range("E2").Value = '  +2.08%  '

$ws.Range("D3").Value = '2.985.07'
$ws.Range("E3").Value = '  +0.75%  '

$ws.Range("E4").Value = '  +0.31%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '562.62'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.91%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '138.49'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +6.77%  '

$ws.Range("E7").Value = '  -0.10%  '

$ws.Range("E8").Value = '  +1.81%  '

$ws.Range("D9").Value = '2.973.44'
$ws.Range("E9").Value = '  +0.53%  '

$ws.Range("E10").Value = '  +2.98%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.16'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +6.62%  '

$ws.Range("E12").Value = '  +1.44%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000229'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.24%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.69'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.46%  '

$ws.Range("E15").Value = '  +2.24%  '

$ws.Range("D16").Value = '3.477.59'
$ws.Range("E16").Value = '  +1.31%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.25'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +8.26%  '

$ws.Range("D18").Value = '2.981.27'
$ws.Range("E18").Value = '  +0.91%  '

$ws.Range("D19").Value = '59.151.40'
$ws.Range("E19").Value = '  +2.13%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '427.88'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.15%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.61'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.43%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.716'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +4.81%  '

$ws.Range("E23").Value = '  +1.56%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.38'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.38%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '80.89'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.69%  '

$ws.Range("E26").Value = '  -0.21%  '

$ws.Range("E27").Value = '  +0.56%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.16'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +9.26%  '

$ws.Range("E29").Value = '  +1.36%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.71'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.81%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '25.70'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.82%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.08'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.19%  '

$ws.Range("E33").Value = '  -5.01%  '

$ws.Range("D34").Value = '0.0₃0769'
$ws.Range("E34").Value = '  +15.77%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.990'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +5.34%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.87'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.14%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.07'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.88%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '49.13'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.45%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.64'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.78%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.71'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +5.32%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '400.71'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +6.56%  '

$ws.Range("D42").Value = '2.768.22'
$ws.Range("E42").Value = '  +4.30%  '

$ws.Range("E43").Value = '  +1.18%  '

$ws.Range("E44").Value = '  -0.39%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.251'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +5.65%  '

$ws.Range("E46").Value = '  -0.04%  '

$ws.Range("B47").Value = 'Stellar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.110'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.11%  '

$ws.Range("B48").Value = 'Arweave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '34.02'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +18.21%  '

$ws.Range("B49").Value = 'Monero'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '120.80'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.36%  '

$ws.Range("E50").Value = '  +0.23%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '23.42'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.09%  '
